# feat: add 2022-Q1 data
#
# 1. Create a new sheet "2022-Q1" (positioned right before "总计") by
#    copying the "2021-Q4" sheet so it inherits identical sheet/cell
#    formatting (sheetPr, styles, margins, etc.), then overwrite its
#    contents with the 2022-Q1 fund holdings (only 4 data rows).
# 2. Insert a new first data row into the "总计" (grand total) sheet for
#    2022-Q1 and renumber the existing index column.

$wb = $excel.ActiveWorkbook

# Helper: write $text into $cell as a genuine text value (never inferred
# as a number/date), and without leaving behind any extra number-format /
# quote-prefix style on the cell (which a plain `$cell.Value = $text`
# assignment would do for numeric-looking strings like "006323").
# Trick: put a `="<text>"` formula in a scratch cell far outside the used
# range, copy it, and PasteSpecial **values only** into the destination -
# that lands a plain text cell with no style baggage, then wipe the
# scratch cell.
function Set-TextValue($sheet, $cell, $text) {
    $scratch = $sheet.Cells.Item(2000, 2000)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

$template = $wb.Worksheets.Item("2021-Q4")
$totalBeforeCopy = $wb.Worksheets.Item("总计")

$template.Copy($totalBeforeCopy)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# NOTE: worksheet references resolve by position, not by identity, so
# after the sheet copy/insert above shifted everything, we must re-fetch
# the "总计" sheet by name rather than reuse $totalBeforeCopy (which now
# points at whatever sheet occupies that original index).
$total = $wb.Worksheets.Item("总计")

# Template sheet has 7 data rows (rows 2-8); 2022-Q1 only needs 4 (rows 2-5).
$newSheet.Rows("6:8").Delete()

# Header row
Set-TextValue $newSheet $newSheet.Cells.Item(1,2) "基金代码"
Set-TextValue $newSheet $newSheet.Cells.Item(1,3) "基金名称"
Set-TextValue $newSheet $newSheet.Cells.Item(1,4) "基金规模"
Set-TextValue $newSheet $newSheet.Cells.Item(1,5) "股票总仓位"
Set-TextValue $newSheet $newSheet.Cells.Item(1,6) "仓位占比"
Set-TextValue $newSheet $newSheet.Cells.Item(1,7) "持有市值(亿元)"
Set-TextValue $newSheet $newSheet.Cells.Item(1,8) "仓位排名"

$fundRows = @(
    @{ code = "006323"; name = "合煦智远嘉选混合A"; size = "1.54"; pos = "79.45"; pct = "3.50"; mv = "0.0539"; rank = 5 },
    @{ code = "501007"; name = "汇添富中证互联网医疗主题指数（LOF）A"; size = "0.58"; pos = "93.89"; pct = "5.19"; mv = "0.0301"; rank = 4 },
    @{ code = "006324"; name = "合煦智远嘉选混合C"; size = "0.59"; pos = "79.45"; pct = "3.50"; mv = "0.0206"; rank = 5 },
    @{ code = "501008"; name = "汇添富中证互联网医疗主题指数（LOF）C"; size = "0.19"; pos = "93.89"; pct = "5.19"; mv = "0.0099"; rank = 4 }
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r,1).Value = $r - 2
    Set-TextValue $newSheet $newSheet.Cells.Item($r,2) $row.code
    Set-TextValue $newSheet $newSheet.Cells.Item($r,3) $row.name
    Set-TextValue $newSheet $newSheet.Cells.Item($r,4) $row.size
    Set-TextValue $newSheet $newSheet.Cells.Item($r,5) $row.pos
    Set-TextValue $newSheet $newSheet.Cells.Item($r,6) $row.pct
    Set-TextValue $newSheet $newSheet.Cells.Item($r,7) $row.mv
    $newSheet.Cells.Item($r,8).Value = $row.rank
    $r = $r + 1
}

# Now update the "总计" (grand total) sheet: insert a new row 2 for
# 2022-Q1 and push the other quarters down.
$total.Rows(2).Insert()
$total.Rows(2).ClearFormats()

$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
Set-TextValue $total $total.Cells.Item(2,2) "2022-Q1"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 0.11

for ($row = 3; $row -le 7; $row++) {
    $total.Cells.Item($row,1).Value = $row - 2
}

# Restore the originally active sheet/selection (copying a sheet makes the
# new copy active; put focus back on the first sheet like the source file).
$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(1).Range("A1").Select() | Out-Null
